$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must remain text;
# force text format before assignment so Excel does not coerce them to numbers.
$textForceCells = @("D5","D6","D7","D16","D17","D19","D20","D22","D23","D24","D25","D26","D29","D31","D32","D33","D35","D40","D41","D43","D44","D45","D46","D49","D51")
foreach ($c in $textForceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "57.006.34"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.414.20"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "489.75"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "154.29"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +18.40%  "
$ws.Range("D9").Value = "2.435.85"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("E10").Value = "  +9.92%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "2.839.63"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "57.004.54"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "20.77"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "2.429.89"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "324.30"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "58.27"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "2.532.86"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  -3.44%  "
$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "150.48"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "18.61"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("E39").Value = "  +9.47%  "
$ws.Range("D40").Value = "34.25"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "0.598"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "268.79"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "0.0533"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "4.58"
$ws.Range("E49").Value = "  -5.52%  "
$ws.Range("D50").Value = "1.878.55"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "17.48"
$ws.Range("E51").Value = "  -1.34%  "

# Restore default (General) style on the cells we temporarily formatted as text
foreach ($c in $textForceCells) { $ws.Range($c).Style = "Normal" }
